$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.907.30'
$ws.Range('E2').Value = '  -0.77%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.869.63'
$ws.Range('E3').Value = '  +0.21%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '305.61'
$ws.Range('E5').Value = '  -0.16%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.0000'
$ws.Range('E6').Value = '  -0.08%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5083'
$ws.Range('E7').Value = '  -1.00%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3660'
$ws.Range('E8').Value = '  -2.68%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07174'
$ws.Range('E9').Value = '  +0.54%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8893'
$ws.Range('E10').Value = '  +0.07%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.57'
$ws.Range('E11').Value = '  -0.55%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.877.87'
$ws.Range('E12').Value = '  +0.73%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.07486'
$ws.Range('E13').Value = '  -0.80%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '94.39'
$ws.Range('E14').Value = '  +5.44%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.218'
$ws.Range('E15').Value = '  -1.70%  '

$ws.Range('E16').Value = '  -0.09%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008488'
$ws.Range('E17').Value = '  +0.22%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '14.12'
$ws.Range('E18').Value = '  +0.22%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.9999'
$ws.Range('E19').Value = '  -0.10%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '26.946.50'
$ws.Range('E20').Value = '  -0.76%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.002'
$ws.Range('E21').Value = '  +0.06%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.116.30'
$ws.Range('E22').Value = '  +1.17%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.34'
$ws.Range('E23').Value = '  -1.25%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.374'
$ws.Range('E24').Value = '  -1.11%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '147.79'
$ws.Range('E25').Value = '  +1.76%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.777'
$ws.Range('E26').Value = '  -3.30%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.84'
$ws.Range('E27').Value = '  -0.62%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.080'
$ws.Range('E28').Value = '  -0.34%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '113.20'
$ws.Range('E29').Value = '  +0.33%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.681'
$ws.Range('E30').Value = '  +0.67%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.701'
$ws.Range('E31').Value = '  +0.73%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09114'
$ws.Range('E32').Value = '  -0.88%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05030'
$ws.Range('E33').Value = '  -1.39%  '

$ws.Range('B34').Value = 'HuobiToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.985'
$ws.Range('E34').Value = '  -3.01%  '

$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7476'
$ws.Range('E35').Value = '  +3.46%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.152'
$ws.Range('E36').Value = '  -0.34%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.221'
$ws.Range('E37').Value = '  +4.15%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5628'
$ws.Range('E38').Value = '  +6.77%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.501'
$ws.Range('E39').Value = '  +0.37%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01987'
$ws.Range('E40').Value = '  -2.45%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.071'
$ws.Range('E41').Value = '  -0.38%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.587'
$ws.Range('E42').Value = '  +1.54%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '115.01'
$ws.Range('E43').Value = '  -1.29%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.543'
$ws.Range('E44').Value = '  +2.96%  '

$ws.Range('E45').Value = '  +1.17%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4768'
$ws.Range('E46').Value = '  +3.30%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.9997'
$ws.Range('E47').Value = '  -0.08%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '10.08'
$ws.Range('E48').Value = '  +1.24%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.553'
$ws.Range('E49').Value = '  -0.45%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '36.88'
$ws.Range('E50').Value = '  +0.84%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '62.92'
$ws.Range('E51').Value = '  -0.83%  '
